$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-03 Monday", "2025-03-04 Tuesday"),
    @("27×38=1026", "40×53=2120"),
    @("97×33=3201", "30×32=960"),
    @("72×22=1584", "48×17=816"),
    @("34×82=2788", "59×81=4779"),
    @("76×48=3648", "63×62=3906"),
    @("90×65=5850", "29×83=2407"),
    @("85×64=5440", "12×75=900"),
    @("66×23=1518", "12×17=204"),
    @("91×25=2275", "16×16=256"),
    @("49×39=1911", "46×38=1748"),
    @("84×36=3024", "94×25=2350"),
    @("27×46=1242", "89×11=979"),
    @("44×87=3828", "92×79=7268"),
    @("60×71=4260", "21×21=441"),
    @("57×59=3363", "28×41=1148"),
    @("33×90=2970", "82×93=7626"),
    @("96×84=8064", "71×69=4899"),
    @("46×29=1334", "55×33=1815"),
    @("49×98=4802", "61×30=1830"),
    @("50×28=1400", "88×87=7656"),
    @("29×18=522", "69×62=4278"),
    @("80×94=7520", "80×45=3600"),
    @("61×67=4087", "27×12=324"),
    @("15×41=615", "25×49=1225"),
    @("15×40=600", "19×32=608")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
